# Generate Report for Handoff
# - Refresh the "Latest Handoff/HO Xliff Generate" timestamps for the rows
#   that just finished handoff generation (0cdc52fa, 1bb760d0, 35bbf8e4,
#   7858f7ed, 7ecafe58, a336ecec -> rows 8,9,10,11,12,14).
# - Mark those same rows' Priority column as "ht" (handoff type) on the
#   per-locale sheets now that a handoff file exists for them.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 11, 12, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-24 08:22:34"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-24 08:22:28"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet: "Latest Handoff Datetime" (H) + "Priority" (E) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-24 08:22:34"
    $wsDeDe.Range("E$r").Value = "ht"
}
